$d = $word.ActiveDocument

# Locate the "ROVER DESCRIPTION" heading paragraph, then the (currently
# empty) paragraph right after it, which is where the new descriptive
# text needs to be added.
$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text
    if ($paraText.Contains("ROVER DESCRIPTION")) {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -eq -1) {
    Write-Host "Could not find ROVER DESCRIPTION heading paragraph"
} else {
    $descParaIndex = $headingIndex + 1
    $descPara = $d.Paragraphs.Item($descParaIndex)

    $bodyText = "Our rover consists of a Terasic DE1-SoC development board with a Cyclone V FPGA/SOC system. An embedded Linux system runs on the embedded ARM platform while interfacing with our custom bitstream loaded onto the FPGA fabric. The FPGA provides a motor controller interface to allow us to control our motors using an original high-level interface. "

    # Add the new run of text to the (empty) paragraph right after the
    # heading, keeping the paragraph's own formatting (Courier New) and
    # giving the run itself the same explicit Courier New formatting.
    $descPara.Range.InsertBefore($bodyText)
    $descPara = $d.Paragraphs.Item($descParaIndex)
    $descPara.Range.Font.Name = "Courier New"
    $descPara.Range.Font.NameBi = "Courier New"

    # Split off a new, blank paragraph (containing a single space) right
    # after the paragraph we just filled in, re-using the same
    # paragraph-level formatting.
    $splitPos = $descPara.Range.End - 1
    $splitRange = $d.Range($splitPos, $splitPos)
    $splitRange.InsertAfter([char]13)

    $spacerPara = $d.Paragraphs.Item($descParaIndex + 1)
    $spacerPara.Range.InsertBefore(" ")
    $spacerPara.Range.Font.Name = "Courier New"
    $spacerPara.Range.Font.NameBi = "Courier New"

    Write-Host "Inserted rover description text after paragraph $headingIndex"
}
